$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Balans")

# --- Fix bugs in DiceItemGenerator (Balans sheet) ---

# D2 referenced B9 (circular-ish / wrong cell); should reference B8
$ws.Range("D2").Formula = "=ROUND(D1+B8,0)"

# B8 base value corrected from 12 to 14
$ws.Range("B8").Value = 14

# B9 should be derived from the dice table (F4) instead of being a hardcoded literal
$ws.Range("B9").Formula = "=F4"

# Remove stray leftover value in H9 (column H is no longer used)
$ws.Range("H9").ClearContents() | Out-Null

# Extend the dice-roll table (columns E:G) down to rolls 11 and 12
$ws.Range("E10").Value = 11
$ws.Range("F10").Formula = "=ROUND(`$B`$8+(`$D`$2-`$B`$8)/(`$B`$2-`$B`$1)*(E10-`$B`$1),0)"
$ws.Range("F10").NumberFormat = "0.00"
$ws.Range("G10").Formula = "=F10-F9"
$ws.Range("G10").NumberFormat = "0.00"

$ws.Range("E11").Value = 12
$ws.Range("F11").Formula = "=ROUND(`$B`$8+(`$D`$2-`$B`$8)/(`$B`$2-`$B`$1)*(E11-`$B`$1),0)"
$ws.Range("F11").NumberFormat = "0.00"
$ws.Range("G11").Formula = "=F11-F10"
$ws.Range("G11").NumberFormat = "0.00"

# Restore the previous cell selection
$ws.Range("B6").Select() | Out-Null
